$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 31850893.47
$ws.Range("P2").Value = 431566376.52
$ws.Range("Q2").Value = 395852398.32
$ws.Range("R2").Value = 40.7189985411
$ws.Range("S2").Value = 302993311.11
$ws.Range("T2").Value = 302993311.11
$ws.Range("U2").Value = 52.0261735341
$ws.Range("V2").Value = 14109951.37
$ws.Range("W2").Value = 45323583.34
$ws.Range("X2").Value = 2160404.09
$ws.Range("Y2").Value = 35050988.76
$ws.Range("Z2").Value = 34957666.87
$ws.Range("AA2").Value = 3106773.4

$ws.Range("AG2").Value = 1908096.62

$ws.Range("AP2").Value = 47.1500655351
$ws.Range("AQ2").Value = 277.075621552038
$ws.Range("AR2").Value = 177.610700559113
$ws.Range("AS2").Value = 29450360.34
$ws.Range("AT2").Value = 117.295573209713
